# Update vm_pu.xlsx values for Case_2_18 (380 kV case) per commit "case with 380 kV done"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"=1.02; "C"=1.034902331059937; "D"=1.045138032827721; "E"=0.992614727750844; "F"=1.053680773777689; "I"=1.042693233828063; "J"=1.040019034977926; "K"=1.047907115532507; "L"=0.9955398523335997; "M"=1.056426072582815; "N"=1.017319576930855 }
    3 = @{ "B"=1.02; "C"=1.035657988304486; "D"=1.045745396744404; "E"=0.9936372048519299; "F"=1.054418524550014; "I"=1.042885083553713; "J"=1.040418997080361; "K"=1.048326390931843; "L"=0.9963617723202687; "M"=1.056977109150555; "N"=1.017451708051472 }
    4 = @{ "B"=1.02; "C"=1.036147432989547; "D"=1.046138776229722; "E"=0.9942998659930998; "F"=1.054896624278156; "I"=1.04300813748756; "J"=1.040677584110099; "K"=1.048597368528881; "L"=0.9968940712668347; "M"=1.057333722234377; "N"=1.017537127954807 }
    5 = @{ "B"=1.02; "C"=1.036353309499031; "D"=1.046304240650694; "E"=0.994578699834602; "F"=1.055097788977462; "I"=1.043059608543607; "J"=1.040786241275399; "K"=1.048711209066991; "L"=0.9971179600053012; "M"=1.057483653930448; "N"=1.017573019416194 }
    6 = @{ "B"=1.02; "C"=1.036387883698901; "D"=1.046332027940154; "E"=0.994625531979634; "F"=1.055131575414291; "I"=1.043068235439564; "J"=1.040804482156121; "K"=1.048730318752187; "L"=0.9971555583673455; "M"=1.057508828748652; "N"=1.017579044619402 }
    7 = @{ "B"=1.02; "C"=1.03615018347611; "D"=1.046140986831328; "E"=0.994303590798249; "F"=1.05489931157925; "I"=1.0430088262723; "J"=1.040679036201653; "K"=1.048598889981772; "L"=0.9968970624462089; "M"=1.057335725585834; "N"=1.017537607613751 }
    8 = @{ "B"=1.02; "C"=1.035157607630605; "D"=1.045343215609467; "E"=0.9929600610674297; "F"=1.05392994856519; "I"=1.042758294590441; "J"=1.040154247817227; "K"=1.048048877341657; "L"=0.9958175282591056; "M"=1.056612285411805; "N"=1.017364247094521 }
    9 = @{ "B"=1.02; "C"=1.033412351864272; "D"=1.043940394541465; "E"=0.9906006454969559; "F"=1.052227459553298; "I"=1.042308550941839; "J"=1.039227913097261; "K"=1.047077285457167; "L"=0.9939188001724441; "M"=1.055337994356494; "N"=1.017058187909476 }
    10 = @{ "B"=1.02; "C"=1.032251508556775; "D"=1.043007282825508; "E"=0.989033133672735; "F"=1.051096395463353; "I"=1.042003212422738; "J"=1.038609364964743; "K"=1.046428030616491; "L"=0.9926553831429383; "M"=1.054488909535334; "N"=1.016853788040943 }
    11 = @{ "B"=1.02; "C"=1.03174950539222; "D"=1.042603758430172; "E"=0.988355674866747; "F"=1.050607591164511; "I"=1.041869702119529; "J"=1.038341308006237; "K"=1.046146552165097; "L"=0.9921088820399291; "M"=1.05412137347366; "N"=1.016765200801441 }
    12 = @{ "B"=1.02; "C"=1.031563138585627; "D"=1.042453951497176; "E"=0.9881042295826724; "F"=1.050426173003965; "I"=1.041819916537782; "J"=1.038241707756098; "K"=1.046041947898131; "L"=0.9919059725120875; "M"=1.053984874521004; "N"=1.016732283879944 }
    13 = @{ "B"=1.02; "C"=1.031603110363974; "D"=1.042486081940856; "E"=0.9881581567098651; "F"=1.050465081194789; "I"=1.041830604474391; "J"=1.038263073772916; "K"=1.046064388136369; "L"=0.9919494934313052; "M"=1.054014153062215; "N"=1.016739345192428 }
    14 = @{ "B"=1.02; "C"=1.031734098212571; "D"=1.042591373698947; "E"=0.9883348863814464; "F"=1.050592592101464; "I"=1.04186559077829; "J"=1.038333075662832; "K"=1.046137906568258; "L"=0.9920921077337197; "M"=1.054110090003117; "N"=1.016762480114008 }
    15 = @{ "B"=1.02; "C"=1.03181481739403; "D"=1.042656258085948; "E"=0.9884438009545853; "F"=1.050671175128163; "I"=1.041887121326869; "J"=1.03837620194202; "K"=1.04618319703452; "L"=0.9921799884222134; "M"=1.054169202662717; "N"=1.016776732769102 }
    16 = @{ "B"=1.02; "C"=1.032284838704569; "D"=1.043034074519072; "E"=0.9890781214508737; "F"=1.051128856081046; "I"=1.04201204581895; "J"=1.0386271504677; "K"=1.046446704232575; "L"=0.9926916645766087; "M"=1.054513304442898; "N"=1.016859665621872 }
    17 = @{ "B"=1.02; "C"=1.032579845938194; "D"=1.043271209186462; "E"=0.989476357848556; "F"=1.051416204199178; "I"=1.042090061125424; "J"=1.03878450547082; "K"=1.046611903519754; "L"=0.9930127773699352; "M"=1.054729184662426; "N"=1.016911665912389 }
    18 = @{ "B"=1.02; "C"=1.03275198121943; "D"=1.043409575781095; "E"=0.9897087662937556; "F"=1.051583901386344; "I"=1.042135440959381; "J"=1.038876266518352; "K"=1.046708227907452; "L"=0.9932001317071769; "M"=1.054855115716164; "N"=1.01694198898063 }
    19 = @{ "B"=1.02; "C"=1.032810685453669; "D"=1.043456763626491; "E"=0.9897880325774034; "F"=1.05164109731272; "I"=1.042150893050751; "J"=1.038907550983386; "K"=1.046741066297075; "L"=0.9932640239640975; "M"=1.054898056893868; "N"=1.016952327022774 }
    20 = @{ "B"=1.02; "C"=1.032548187969961; "D"=1.043245761699346; "E"=0.9894336180360679; "F"=1.051385364944873; "I"=1.042081703762214; "J"=1.038767624978712; "K"=1.046594182660504; "L"=0.9929783193494215; "M"=1.054706021527211; "N"=1.016906087579068 }
    21 = @{ "B"=1.02; "C"=1.031695522793291; "D"=1.042560365676243; "E"=0.9882828385668249; "F"=1.050555039282821; "I"=1.041855293529412; "J"=1.038312462722657; "K"=1.04611625860852; "L"=0.9920501090198102; "M"=1.05408183839034; "N"=1.016755667774901 }
    22 = @{ "B"=1.02; "C"=1.031159995098409; "D"=1.042129893802619; "E"=0.9875604150241495; "F"=1.05003382360636; "I"=1.0417118189714; "J"=1.038026099886778; "K"=1.045815476541974; "L"=0.9914670000341481; "M"=1.053689508039782; "N"=1.016661025542497 }
    23 = @{ "B"=1.02; "C"=1.031443833152843; "D"=1.042358050457186; "E"=0.9879432794643023; "F"=1.050310049254104; "I"=1.041787983548128; "J"=1.03817792326875; "K"=1.045974954018165; "L"=0.991776070289318; "M"=1.053897477967012; "N"=1.016711203411091 }
    24 = @{ "B"=1.02; "C"=1.032562492650072; "D"=1.043257260169312; "E"=0.9894529299347244; "F"=1.051399299592797; "I"=1.042085480481792; "J"=1.038775252612973; "K"=1.046602190057904; "L"=0.9929938892766442; "M"=1.054716487914823; "N"=1.016908608212005 }
    25 = @{ "B"=1.02; "C"=1.033863081616542; "D"=1.044302695368733; "E"=0.9912096547607049; "F"=1.052666910703025; "I"=1.042425795490298; "J"=1.039467573297457; "K"=1.047328740686164; "L"=0.9944092447426414; "M"=1.055667358769473; "N"=1.017137377023886 }
}

foreach ($rowKey in $data.Keys) {
    $rowNum = [int]$rowKey
    $rowData = $data[$rowKey]
    foreach ($col in $rowData.Keys) {
        $ws.Range("$col$rowNum").Value = $rowData[$col]
    }
}

$wb.Save()